# Update: po 14. 02. 2022
# Revises AgTests (F) / AgPosit (G) figures for several existing rows and
# appends three new daily rows (709-711) for 2022-02-11, 2022-02-12 and
# 2022-02-13 (serials 44603-44605).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised values on existing rows ---
$ws.Cells.Item(603, 6).Value = 32362
$ws.Cells.Item(604, 6).Value = 30349
$ws.Cells.Item(614, 6).Value = 47781
$ws.Cells.Item(615, 6).Value = 37062
$ws.Cells.Item(616, 6).Value = 38477
$ws.Cells.Item(617, 6).Value = 39172
$ws.Cells.Item(621, 6).Value = 56512
$ws.Cells.Item(630, 6).Value = 46881
$ws.Cells.Item(645, 6).Value = 35749
$ws.Cells.Item(649, 6).Value = 62674
$ws.Cells.Item(651, 6).Value = 37137
$ws.Cells.Item(653, 6).Value = 34156
$ws.Cells.Item(657, 6).Value = 34086
$ws.Cells.Item(680, 6).Value = 28473
$ws.Cells.Item(686, 6).Value = 34428
$ws.Cells.Item(692, 6).Value = 41528
$ws.Cells.Item(693, 6).Value = 39434
$ws.Cells.Item(694, 6).Value = 37469
$ws.Cells.Item(694, 7).Value = 2773
$ws.Cells.Item(695, 6).Value = 36850
$ws.Cells.Item(695, 7).Value = 3109
$ws.Cells.Item(696, 6).Value = 17674
$ws.Cells.Item(696, 7).Value = 2193
$ws.Cells.Item(697, 6).Value = 28583
$ws.Cells.Item(697, 7).Value = 2998
$ws.Cells.Item(698, 6).Value = 68315
$ws.Cells.Item(698, 7).Value = 5741
$ws.Cells.Item(699, 6).Value = 42841
$ws.Cells.Item(699, 7).Value = 4259
$ws.Cells.Item(700, 6).Value = 42964
$ws.Cells.Item(700, 7).Value = 4208
$ws.Cells.Item(701, 6).Value = 41227
$ws.Cells.Item(701, 7).Value = 3791
$ws.Cells.Item(702, 6).Value = 35445
$ws.Cells.Item(702, 7).Value = 3826
$ws.Cells.Item(703, 6).Value = 16522
$ws.Cells.Item(703, 7).Value = 2520
$ws.Cells.Item(704, 6).Value = 24124
$ws.Cells.Item(704, 7).Value = 3562
$ws.Cells.Item(705, 6).Value = 53694
$ws.Cells.Item(705, 7).Value = 6084
$ws.Cells.Item(706, 6).Value = 39354
$ws.Cells.Item(706, 7).Value = 4776
$ws.Cells.Item(707, 6).Value = 35287
$ws.Cells.Item(707, 7).Value = 3999
$ws.Cells.Item(708, 6).Value = 32060
$ws.Cells.Item(708, 7).Value = 3503

# --- New rows 709-711 ---
$ws.Cells.Item(709, 1).Value = 44603
$ws.Cells.Item(709, 2).Value = 1226925
$ws.Cells.Item(709, 3).Value = 35680
$ws.Cells.Item(709, 4).Value = 19872
$ws.Cells.Item(709, 5).Value = 18081
$ws.Cells.Item(709, 6).Value = 24996
$ws.Cells.Item(709, 7).Value = 3034

$ws.Cells.Item(710, 1).Value = 44604
$ws.Cells.Item(710, 2).Value = 1241121
$ws.Cells.Item(710, 3).Value = 24925
$ws.Cells.Item(710, 4).Value = 14196
$ws.Cells.Item(710, 5).Value = 18095
$ws.Cells.Item(710, 6).Value = 10547
$ws.Cells.Item(710, 7).Value = 1939

$ws.Cells.Item(711, 1).Value = 44605
$ws.Cells.Item(711, 2).Value = 1249673
$ws.Cells.Item(711, 3).Value = 15645
$ws.Cells.Item(711, 4).Value = 8552
$ws.Cells.Item(711, 5).Value = 18105
$ws.Cells.Item(711, 6).Value = 12299
$ws.Cells.Item(711, 7).Value = 2265
